$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BIS-769: the SAMPLE_TYPE property-types table (row 4 = header row) gains
# two new trailing columns, "Pattern" (M) and "Pattern Type" (N), matching
# the bold header style already used by "Unique" (L4).
$ws.Range("M4").Value = "Pattern"
$ws.Range("N4").Value = "Pattern Type"

$ws.Range("M4:N4").Font.Size = 11
$ws.Range("M4:N4").Font.Bold = $true

# The active selection moves to the newly added header cells.
$ws.Range("M4:N4").Select()
